# Update the "Class Name" column (K) for the student import sample rows.
# Previously every row used the placeholder "SE1902"; split into distinct
# program codes per row: SE (row 2), BE (row 3), SAP (row 4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "SE"
$ws.Range("K3").Value = "BE"
$ws.Range("K4").Value = "SAP"

# Match the author's final selection/active cell in the saved file.
$null = $ws.Range("K14").Select()
